$d = $word.ActiveDocument

# Paragraph 8 is the originally-empty paragraph that becomes the first new
# figure legend paragraph; a second new paragraph is created right after it via a
# paragraph-mark insertion, which in turn pushes the existing blank paragraph down.
$targetPara = $d.Paragraphs.Item(8)
$targetRange = $targetPara.Range

$markerText = "@@M001X@@@@M002X@@@@M003X@@@@M004X@@@@M005X@@@@M006X@@@@M007X@@@@M008X@@@@M009X@@@@M010X@@@@M011X@@@@M012X@@@@M013X@@@@M014X@@@@M015X@@@@M016X@@@@M017X@@@@M018X@@@@M019X@@@@M020X@@@@M021X@@@@M022X@@@@M023X@@@@M024X@@@@M025X@@@@M026X@@@@M027X@@@@M028X@@@@M029X@@@@M030X@@@@M031X@@@@M032X@@@@M033X@@@@M034X@@@@M035X@@@@M036X@@@@M037X@@@@M038X@@@@M039X@@@@M040X@@@@M041X@@@@M042X@@@@M043X@@@@M044X@@@@M045X@@@@M046X@@@@M047X@@@@M048X@@@@M049X@@@@M050X@@@@M051X@@@@M052X@@@@M053X@@@@M054X@@@@M055X@@@@M056X@@@@M057X@@" + [char]13 + "@@M058X@@@@M059X@@@@M060X@@@@M061X@@@@M062X@@@@M063X@@@@M064X@@@@M065X@@@@M066X@@@@M067X@@@@M068X@@@@M069X@@@@M070X@@@@M071X@@@@M072X@@@@M073X@@@@M074X@@@@M075X@@@@M076X@@@@M077X@@@@M078X@@@@M079X@@@@M080X@@@@M081X@@@@M082X@@@@M083X@@@@M084X@@@@M085X@@@@M086X@@@@M087X@@@@M088X@@@@M089X@@@@M090X@@@@M091X@@@@M092X@@@@M093X@@@@M094X@@"
$targetRange.InsertAfter($markerText)
Write-Output ("Paragraph count after split: " + $d.Paragraphs.Count)

function Replace-Marker($Marker, $NewText, $Italic, $Bold) {
    $f = $d.Content.Find
    $f.ClearFormatting()
    $f.Replacement.ClearFormatting()
    $f.Replacement.Font.NameAscii = "Times New Roman"
    $f.Replacement.Font.NameOther = "Times New Roman"
    $f.Replacement.Font.NameBi = "Times New Roman"
    if ($Italic) { $f.Replacement.Font.Italic = $true }
    if ($Bold) { $f.Replacement.Font.Bold = $true }
    $ok = $f.Execute($Marker, $true, $false, $false, $false, $false, $true, 1, $false, $NewText, 2)
    if (-not $ok) {
        Write-Output ("WARNING: replace failed for marker " + $Marker)
    }
}

Replace-Marker "@@M001X@@" "Figure." $false $false
Replace-Marker "@@M002X@@" " " $false $false
Replace-Marker "@@M003X@@" "Phylogenetic" $false $false
Replace-Marker "@@M004X@@" " tree of the amino acid sequence of " $false $false
Replace-Marker "@@M005X@@" "dddD" $false $false
Replace-Marker "@@M006X@@" " " $false $false
Replace-Marker "@@M007X@@" "DMSP " $false $false
Replace-Marker "@@M008X@@" "lyase" $false $false
Replace-Marker "@@M009X@@" " " $false $false
Replace-Marker "@@M010X@@" "homologs" $false $false
Replace-Marker "@@M011X@@" " from Organic Lake and " $false $false
Replace-Marker "@@M012X@@" "public databases with " $false $false
Replace-Marker "@@M013X@@" "E. coli" $true $false
Replace-Marker "@@M014X@@" " " $false $false
Replace-Marker "@@M015X@@" "carnitine" $false $false
Replace-Marker "@@M016X@@" " " $false $false
Replace-Marker "@@M017X@@" "CoA" $false $false
Replace-Marker "@@M018X@@" " " $false $false
Replace-Marker "@@M019X@@" "transferase" $false $false
Replace-Marker "@@M020X@@" " of as an " $false $false
Replace-Marker "@@M021X@@" "outgroup" $false $false
Replace-Marker "@@M022X@@" ". The tree was computed from a 75 residue region within the conserved amino-terminal class III coenzyme A domain (" $false $false
Replace-Marker "@@M023X@@" "CaiB" $false $false
Replace-Marker "@@M024X@@" ") using the neighbor-joining algorithm. Bootstrap values are shown at the nodes. Organic Lake sequences from this study are marked with an asterisk (*)." $false $false
Replace-Marker "@@M025X@@" " Numbers in parentheses are counts of sequences which clustered with the Organic Lake homolog shown in the tree with 90% amino acid identity." $false $false
Replace-Marker "@@M026X@@" " Sequences with confirmed DMSP " $false $false
Replace-Marker "@@M027X@@" "lyase" $false $false
Replace-Marker "@@M028X@@" " activity are shown in bold. " $false $false
Replace-Marker "@@M029X@@" "Accession numbers from top to bottom are: " $false $false
Replace-Marker "@@M030X@@" "EBA01716.1" $false $false
Replace-Marker "@@M031X@@" ", " $false $false
Replace-Marker "@@M032X@@" "AEV37420.1" $false $false
Replace-Marker "@@M033X@@" ", " $false $false
Replace-Marker "@@M034X@@" "ACY01992.1" $false $false
Replace-Marker "@@M035X@@" ", " $false $false
Replace-Marker "@@M036X@@" "ADZ91595.1" $false $false
Replace-Marker "@@M037X@@" ", " $false $false
Replace-Marker "@@M038X@@" "EAQ63474.1" $false $false
Replace-Marker "@@M039X@@" ", " $false $false
Replace-Marker "@@M040X@@" "ABR72937.1" $false $false
Replace-Marker "@@M041X@@" ", " $false $false
Replace-Marker "@@M042X@@" "ACV84065.1" $false $false
Replace-Marker "@@M043X@@" ", " $false $false
Replace-Marker "@@M044X@@" "ACY02894.1" $false $false
Replace-Marker "@@M045X@@" ", " $false $false
Replace-Marker "@@M046X@@" "ABI89851.1" $false $false
Replace-Marker "@@M047X@@" ", " $false $false
Replace-Marker "@@M048X@@" "YP_002822700.1" $false $false
Replace-Marker "@@M049X@@" ", " $false $false
Replace-Marker "@@M050X@@" "EEE36156.1" $false $false
Replace-Marker "@@M051X@@" ", " $false $false
Replace-Marker "@@M052X@@" "ABV95365.1" $false $false
Replace-Marker "@@M053X@@" ", " $false $false
Replace-Marker "@@M054X@@" "AAV94987.1" $false $false
Replace-Marker "@@M055X@@" " and " $false $false
Replace-Marker "@@M056X@@" "EGB36199.1" $false $false
Replace-Marker "@@M057X@@" "." $false $false
Replace-Marker "@@M058X@@" "Figure." $false $false
Replace-Marker "@@M059X@@" " " $false $false
Replace-Marker "@@M060X@@" "Phylogenetic" $false $false
Replace-Marker "@@M061X@@" " tree of the amino acid sequence of " $false $false
Replace-Marker "@@M062X@@" "dddL" $false $false
Replace-Marker "@@M063X@@" " DMSP " $false $false
Replace-Marker "@@M064X@@" "lyase" $false $false
Replace-Marker "@@M065X@@" " " $false $false
Replace-Marker "@@M066X@@" "homologs" $false $false
Replace-Marker "@@M067X@@" " from Organic Lake and public databases." $false $false
Replace-Marker "@@M068X@@" " The tree was computed from an 84 residue N-terminal region using the neighbor-joining algorithm. Bootstrap values are shown at the nodes. Organic Lake sequences from this study are " $false $false
Replace-Marker "@@M069X@@" "marked with an asterisk (*). Numbers in parentheses are counts of sequences which clustered with the Organic Lake homolog shown in the tree with 90% amino acid identity. Sequences with confirmed DMSP " $false $false
Replace-Marker "@@M070X@@" "lyase" $false $false
Replace-Marker "@@M071X@@" " activity are shown in bold. Accession numbers from top to bottom are: " $false $false
Replace-Marker "@@M072X@@" "EEB86351.1" $false $false
Replace-Marker "@@M073X@@" ", " $false $false
Replace-Marker "@@M074X@@" "ADK55772.1" $false $false
Replace-Marker "@@M075X@@" ", " $false $false
Replace-Marker "@@M076X@@" "EAQ07081.1" $false $false
Replace-Marker "@@M077X@@" ", " $false $false
Replace-Marker "@@M078X@@" "EEE47811.1" $false $false
Replace-Marker "@@M079X@@" ", " $false $false
Replace-Marker "@@M080X@@" "EAV43167.1" $false $false
Replace-Marker "@@M081X@@" "," $false $false
Replace-Marker "@@M082X@@" " " $false $false
Replace-Marker "@@M083X@@" "EAU41122.1" $false $false
Replace-Marker "@@M084X@@" ", " $false $false
Replace-Marker "@@M085X@@" "EAQ10619.1" $false $false
Replace-Marker "@@M086X@@" ", " $false $false
Replace-Marker "@@M087X@@" "ABV95046.1" $false $false
Replace-Marker "@@M088X@@" ", " $false $false
Replace-Marker "@@M089X@@" "EAQ04071.1" $false $false
Replace-Marker "@@M090X@@" ", " $false $false
Replace-Marker "@@M091X@@" "ABA77574.1" $false $false
Replace-Marker "@@M092X@@" " and " $false $false
Replace-Marker "@@M093X@@" "EHJ04839.1" $false $false
Replace-Marker "@@M094X@@" "." $false $false

Write-Output "Done applying figure legend paragraphs."
